$wb = $excel.ActiveWorkbook

# --- constants_evaluated (text cells; keep as shared-string / text type) ---
$ws = $wb.Worksheets.Item("constants_evaluated")
$ws.Range("Z1").Formula = '="7.47838163473776"'
$ws.Range("Z1").Copy()
$ws.Range("C4").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="-6.30859375"'
$ws.Range("Z1").Copy()
$ws.Range("B5").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="25288264440.948"'
$ws.Range("Z1").Copy()
$ws.Range("C5").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# --- correlation_matrix ---
$ws = $wb.Worksheets.Item("correlation_matrix")
$ws.Range("B2").Value = -0.0708318160190816
$ws.Range("A3").Value = -0.0708318160190816

# --- adj_r_squared ---
$ws = $wb.Worksheets.Item("adj_r_squared")
$ws.Range("A2").Value = 0.755877491863509

# --- mol_ext_coefficients_calc ---
$ws = $wb.Worksheets.Item("mol_ext_coefficients_calc")
$ws.Range("F2").Value = 82375689655040416
$ws.Range("E3").Value = 36.3301929247834
$ws.Range("F3").Value = 63221652789401264

# --- equilibrium_concentrations ---
$ws = $wb.Worksheets.Item("equilibrium_concentrations")
$ws.Range("A2").Value = 0.0000000401396896388809
$ws.Range("C2").Value = 0.0000277369337910124
$ws.Range("D2").Value = 0.000000000000000000000000353063777218641
$ws.Range("E2").Value = 0.000000377073480650462
$ws.Range("A3").Value = 0.000000191482654670096
$ws.Range("D3").Value = 0.00000000000000000000000658157127997869
$ws.Range("E3").Value = 0.0000000790443004377557
$ws.Range("A4").Value = 0.000000294582072854087
$ws.Range("D4").Value = 0.000000000000000000000013868423836588
$ws.Range("E4").Value = 0.0000000513799510530944
$ws.Range("A5").Value = 0.000000605603648085603
$ws.Range("C5").Value = 0.00022931938895647
$ws.Range("D5").Value = 0.0000000000000000000000440402688631355
$ws.Range("E5").Value = 0.0000000249926045396322
$ws.Range("A6").Value = 0.00000122883861718864
$ws.Range("D6").Value = 0.000000000000000000000121030242540732
$ws.Range("E6").Value = 0.0000000123170058888527
$ws.Range("A7").Value = 0.0000183743973736412
$ws.Range("B7").Value = 0.0000160735736398244
$ws.Range("C7").Value = 0.000457626426360219
$ws.Range("D7").Value = 0.00000000000000000000266652035390122
$ws.Range("E7").Value = 0.000000000823733817037981
$ws.Range("B8").Value = 0.00000215894214303153
$ws.Range("D8").Value = 0.000000000000000000021078183413176
$ws.Range("E8").Value = 0.000000000107375954477657
$ws.Range("A9").Value = 0.000275009109236228
$ws.Range("B9").Value = 0.00000110905419978861
$ws.Range("D9").Value = 0.0000000000000000000412147983676479
$ws.Range("E9").Value = 0.0000000000550367677870657
$ws.Range("B10").Value = 0.00000074530152828605
$ws.Range("C10").Value = 0.000472954700104318
$ws.Range("D10").Value = 0.0000000000000000000614245857044143
$ws.Range("E10").Value = 0.0000000000369571108791822
$ws.Range("B11").Value = 0.000000324102250960064
$ws.Range("D11").Value = 0.000000000000000000141502916929003
$ws.Range("E11").Value = 0.0000000000160568907482791
$ws.Range("B12").Value = 0.000000226259755604239
$ws.Range("D12").Value = 0.000000000000000000202777440513961
$ws.Range("E12").Value = 0.0000000000112071958982123
$ws.Range("B13").Value = 0.0000000645513357276695
$ws.Range("C13").Value = 0.000473635448718473
$ws.Range("D13").Value = 0.000000000000000000711243517887338
$ws.Range("E13").Value = 0.00000000000319629298020497

# --- absorbance_calc_abs_errors ---
$ws = $wb.Worksheets.Item("absorbance_calc_abs_errors")
$ws.Range("H2").Value = 1.66249649800997
$ws.Range("I2").Value = 1.66894105471416
$ws.Range("K2").Value = 1.67276526565934
$ws.Range("D3").Value = 0.0385940256461668
$ws.Range("E3").Value = 0.0477851934321829
$ws.Range("F3").Value = 0.0344252613773512
$ws.Range("G3").Value = 0.0272113771163245
$ws.Range("H3").Value = 0.00249649800997154
$ws.Range("I3").Value = -0.0240589452858353
$ws.Range("J3").Value = -0.0460283582061767
$ws.Range("K3").Value = -0.0502347343406624
$ws.Range("L3").Value = 0.0195109358960501
$ws.Range("M3").Value = 0.0465931181537631
$ws.Range("N3").Value = -0.00946436872104739

# --- absorbance_calc_rel_errors ---
$ws = $wb.Worksheets.Item("absorbance_calc_rel_errors")
$ws.Range("H2").Value = 1.66249649800997
$ws.Range("I2").Value = 1.66894105471416
$ws.Range("K2").Value = 1.67276526565934
$ws.Range("D3").Value = 0.0257293504307778
$ws.Range("E3").Value = 0.0317509590911514
$ws.Range("F3").Value = 0.0222529162103111
$ws.Range("G3").Value = 0.0171897518106914
$ws.Range("H3").Value = 0.00150391446383828
$ws.Range("I3").Value = -0.0142108359632813
$ws.Range("J3").Value = -0.0268074305219433
$ws.Range("K3").Value = -0.029155388473977
$ws.Range("L3").Value = 0.0117535758409941
$ws.Range("M3").Value = 0.0284451270779995
$ws.Range("N3").Value = -0.00545182530014251
